$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between paired rows; column A (Indice) stays put ---
# swap row 8 and 9
$ws.Range("F8").Value = 'Feyenoord'
$ws.Range("F9").Value = 'AZ Alkmaar'
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 5
$ws.Range("H8").Value = 'Sittard'
$ws.Range("H9").Value = 'G.A. Eagles'
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J8").Value = 1.22
$ws.Range("J9").Value = 1.37
$ws.Range("K8").Value = '07/07/2023 11:12'
$ws.Range("K9").Value = '07/07/2023 11:12'
$ws.Range("L8").Value = 1.19
$ws.Range("L9").Value = 1.54
$ws.Range("M8").Value = '13/08/2023 14:25'
$ws.Range("M9").Value = '13/08/2023 14:27'
$ws.Range("N8").Value = 7.61
$ws.Range("N9").Value = 5.51
$ws.Range("O8").Value = '07/07/2023 11:12'
$ws.Range("O9").Value = '07/07/2023 11:12'
$ws.Range("P8").Value = 7.38
$ws.Range("P9").Value = 4.41
$ws.Range("Q8").Value = '13/08/2023 14:29'
$ws.Range("Q9").Value = '13/08/2023 14:26'
$ws.Range("R8").Value = 12.86
$ws.Range("R9").Value = 7.89
$ws.Range("S8").Value = '07/07/2023 11:12'
$ws.Range("S9").Value = '07/07/2023 11:12'
$ws.Range("T8").Value = 16.76
$ws.Range("T9").Value = 6.35
$ws.Range("U8").Value = '13/08/2023 14:29'
$ws.Range("U9").Value = '13/08/2023 14:26'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/feyenoord-sittard/6ZGlMZuC/'
$ws.Range("V9").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/az-alkmaar-g-a-eagles/SrIpNFQ5/'

# swap row 17 and 18
$ws.Range("F17").Value = 'Twente'
$ws.Range("F18").Value = 'Sparta Rotterdam'
$ws.Range("G17").Value = 3
$ws.Range("G18").Value = 2
$ws.Range("H17").Value = 'Zwolle'
$ws.Range("H18").Value = 'Feyenoord'
$ws.Range("I17").Value = 1
$ws.Range("I18").Value = 2
$ws.Range("J17").Value = 1.35
$ws.Range("J18").Value = 5
$ws.Range("K17").Value = '13/08/2023 17:12'
$ws.Range("K18").Value = '13/08/2023 14:42'
$ws.Range("L17").Value = 1.44
$ws.Range("L18").Value = 4.94
$ws.Range("M17").Value = '20/08/2023 14:21'
$ws.Range("M18").Value = '20/08/2023 14:28'
$ws.Range("N17").Value = 5.56
$ws.Range("N18").Value = 4.52
$ws.Range("O17").Value = '13/08/2023 17:12'
$ws.Range("O18").Value = '13/08/2023 14:42'
$ws.Range("P17").Value = 5.21
$ws.Range("P18").Value = 4.35
$ws.Range("Q17").Value = '20/08/2023 14:28'
$ws.Range("Q18").Value = '20/08/2023 14:29'
$ws.Range("R17").Value = 8.79
$ws.Range("R18").Value = 1.61
$ws.Range("S17").Value = '13/08/2023 17:12'
$ws.Range("S18").Value = '13/08/2023 14:42'
$ws.Range("T17").Value = 6.95
$ws.Range("T18").Value = 1.67
$ws.Range("U17").Value = '20/08/2023 14:29'
$ws.Range("U18").Value = '20/08/2023 14:28'
$ws.Range("V17").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-zwolle/KxsOzZf5/'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-feyenoord/21WSZhAB/'

# swap row 49 and 50
$ws.Range("F49").Value = 'Waalwijk'
$ws.Range("F50").Value = 'Zwolle'
$ws.Range("G49").Value = 1
$ws.Range("G50").Value = 0
$ws.Range("H49").Value = 'Twente'
$ws.Range("H50").Value = 'AZ Alkmaar'
$ws.Range("I49").Value = 0
$ws.Range("I50").Value = 3
$ws.Range("J49").Value = 4.49
$ws.Range("J50").Value = 4.72
$ws.Range("K49").Value = '17/09/2023 13:43'
$ws.Range("K50").Value = '17/09/2023 16:13'
$ws.Range("L49").Value = 6.03
$ws.Range("L50").Value = 4.63
$ws.Range("M49").Value = '24/09/2023 16:43'
$ws.Range("M50").Value = '24/09/2023 16:39'
$ws.Range("N49").Value = 4.43
$ws.Range("N50").Value = 4.32
$ws.Range("O49").Value = '17/09/2023 13:43'
$ws.Range("O50").Value = '17/09/2023 16:13'
$ws.Range("P49").Value = 4.72
$ws.Range("P50").Value = 4.03
$ws.Range("Q49").Value = '24/09/2023 16:44'
$ws.Range("Q50").Value = '24/09/2023 16:39'
$ws.Range("R49").Value = 1.68
$ws.Range("R50").Value = 1.67
$ws.Range("S49").Value = '17/09/2023 13:43'
$ws.Range("S50").Value = '17/09/2023 16:13'
$ws.Range("T49").Value = 1.53
$ws.Range("T50").Value = 1.77
$ws.Range("U49").Value = '24/09/2023 16:38'
$ws.Range("U50").Value = '24/09/2023 16:39'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-twente/OWEoHsa8/'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-az-alkmaar/ATQ3hbM7/'

# swap row 56 and 57
$ws.Range("F56").Value = 'PSV'
$ws.Range("F57").Value = 'Utrecht'
$ws.Range("G56").Value = 3
$ws.Range("G57").Value = 0
$ws.Range("H56").Value = 'FC Volendam'
$ws.Range("H57").Value = 'Almere City'
$ws.Range("I56").Value = 1
$ws.Range("I57").Value = 2
$ws.Range("J56").Value = 1.06
$ws.Range("J57").Value = 1.51
$ws.Range("K56").Value = '27/09/2023 18:12'
$ws.Range("K57").Value = '23/09/2023 19:12'
$ws.Range("L56").Value = 1.04
$ws.Range("L57").Value = 1.57
$ws.Range("M56").Value = '30/09/2023 14:21'
$ws.Range("M57").Value = '30/09/2023 18:44'
$ws.Range("N56").Value = 17.77
$ws.Range("N57").Value = 4.75
$ws.Range("O56").Value = '27/09/2023 18:12'
$ws.Range("O57").Value = '23/09/2023 19:12'
$ws.Range("P56").Value = 25.88
$ws.Range("P57").Value = 4.4
$ws.Range("Q56").Value = '30/09/2023 18:44'
$ws.Range("Q57").Value = '30/09/2023 18:44'
$ws.Range("R56").Value = 24.06
$ws.Range("R57").Value = 5.79
$ws.Range("S56").Value = '27/09/2023 18:12'
$ws.Range("S57").Value = '23/09/2023 19:12'
$ws.Range("T56").Value = 42.53
$ws.Range("T57").Value = 5.93
$ws.Range("U56").Value = '30/09/2023 18:44'
$ws.Range("U57").Value = '30/09/2023 18:44'
$ws.Range("V56").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-fc-volendam/EFCT8J6l/'
$ws.Range("V57").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/utrecht-almere-city/dv2Y7wMf/'

# swap row 60 and 61
$ws.Range("F60").Value = 'Excelsior'
$ws.Range("F61").Value = 'Heracles'
$ws.Range("G60").Value = 2
$ws.Range("G61").Value = 2
$ws.Range("H60").Value = 'Sparta Rotterdam'
$ws.Range("H61").Value = 'Zwolle'
$ws.Range("I60").Value = 1
$ws.Range("I61").Value = 1
$ws.Range("J60").Value = 2.83
$ws.Range("J61").Value = 2.09
$ws.Range("K60").Value = '24/09/2023 13:42'
$ws.Range("K61").Value = '28/09/2023 19:12'
$ws.Range("L60").Value = 3.17
$ws.Range("L61").Value = 2.44
$ws.Range("M60").Value = '01/10/2023 14:20'
$ws.Range("M61").Value = '01/10/2023 14:24'
$ws.Range("N60").Value = 3.52
$ws.Range("N61").Value = 4.05
$ws.Range("O60").Value = '24/09/2023 13:42'
$ws.Range("O61").Value = '28/09/2023 19:12'
$ws.Range("P60").Value = 3.54
$ws.Range("P61").Value = 3.69
$ws.Range("Q60").Value = '01/10/2023 14:20'
$ws.Range("Q61").Value = '01/10/2023 14:27'
$ws.Range("R60").Value = 2.49
$ws.Range("R61").Value = 3.19
$ws.Range("S60").Value = '24/09/2023 13:42'
$ws.Range("S61").Value = '28/09/2023 19:12'
$ws.Range("T60").Value = 2.33
$ws.Range("T61").Value = 2.88
$ws.Range("U60").Value = '01/10/2023 14:20'
$ws.Range("U61").Value = '01/10/2023 14:27'
$ws.Range("V60").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-sparta-rotterdam/KCXlChDt/'
$ws.Range("V61").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/heracles-zwolle/rmALTKrE/'

# swap row 78 and 79
$ws.Range("F78").Value = 'Heracles'
$ws.Range("F79").Value = 'Excelsior'
$ws.Range("G78").Value = 2
$ws.Range("G79").Value = 2
$ws.Range("H78").Value = 'Twente'
$ws.Range("H79").Value = 'Zwolle'
$ws.Range("I78").Value = 2
$ws.Range("I79").Value = 4
$ws.Range("J78").Value = 4.47
$ws.Range("J79").Value = 2.37
$ws.Range("K78").Value = '09/10/2023 14:42'
$ws.Range("K79").Value = '09/10/2023 14:42'
$ws.Range("L78").Value = 6.03
$ws.Range("L79").Value = 2.31
$ws.Range("M78").Value = '22/10/2023 14:28'
$ws.Range("M79").Value = '22/10/2023 14:27'
$ws.Range("N78").Value = 4.05
$ws.Range("N79").Value = 3.56
$ws.Range("O78").Value = '09/10/2023 14:42'
$ws.Range("O79").Value = '09/10/2023 14:42'
$ws.Range("P78").Value = 4.66
$ws.Range("P79").Value = 3.89
$ws.Range("Q78").Value = '22/10/2023 14:28'
$ws.Range("Q79").Value = '22/10/2023 14:28'
$ws.Range("R78").Value = 1.75
$ws.Range("R79").Value = 2.97
$ws.Range("S78").Value = '09/10/2023 14:42'
$ws.Range("S79").Value = '09/10/2023 14:42'
$ws.Range("T78").Value = 1.53
$ws.Range("T79").Value = 2.96
$ws.Range("U78").Value = '22/10/2023 14:22'
$ws.Range("U79").Value = '22/10/2023 14:24'
$ws.Range("V78").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/heracles-twente/SKsRD5lo/'
$ws.Range("V79").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-zwolle/QBeaIotU/'

# swap row 83 and 84
$ws.Range("F83").Value = 'Almere City'
$ws.Range("F84").Value = 'Sparta Rotterdam'
$ws.Range("G83").Value = 0
$ws.Range("G84").Value = 2
$ws.Range("H83").Value = 'G.A. Eagles'
$ws.Range("H84").Value = 'Waalwijk'
$ws.Range("I83").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J83").Value = 2.51
$ws.Range("J84").Value = 1.75
$ws.Range("K83").Value = '22/10/2023 17:12'
$ws.Range("K84").Value = '22/10/2023 17:12'
$ws.Range("L83").Value = 3.13
$ws.Range("L84").Value = 1.78
$ws.Range("M83").Value = '28/10/2023 19:58'
$ws.Range("M84").Value = '28/10/2023 19:35'
$ws.Range("N83").Value = 3.48
$ws.Range("N84").Value = 3.96
$ws.Range("O83").Value = '22/10/2023 17:12'
$ws.Range("O84").Value = '22/10/2023 17:12'
$ws.Range("P83").Value = 3.57
$ws.Range("P84").Value = 3.95
$ws.Range("Q83").Value = '28/10/2023 19:58'
$ws.Range("Q84").Value = '28/10/2023 19:37'
$ws.Range("R83").Value = 2.89
$ws.Range("R84").Value = 4.62
$ws.Range("S83").Value = '22/10/2023 17:12'
$ws.Range("S84").Value = '22/10/2023 17:12'
$ws.Range("T83").Value = 2.34
$ws.Range("T84").Value = 4.63
$ws.Range("U83").Value = '28/10/2023 19:58'
$ws.Range("U84").Value = '28/10/2023 19:35'
$ws.Range("V83").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-g-a-eagles/xjPHzTeb/'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-waalwijk/YHIQY6QA/'

# swap row 87 and 88
$ws.Range("F87").Value = 'FC Volendam'
$ws.Range("F88").Value = 'PSV'
$ws.Range("G87").Value = 3
$ws.Range("G88").Value = 5
$ws.Range("H87").Value = 'Excelsior'
$ws.Range("H88").Value = 'Ajax'
$ws.Range("I87").Value = 1
$ws.Range("I88").Value = 2
$ws.Range("J87").Value = 3.04
$ws.Range("J88").Value = 1.54
$ws.Range("K87").Value = '22/10/2023 14:42'
$ws.Range("K88").Value = '22/10/2023 12:43'
$ws.Range("L87").Value = 2.62
$ws.Range("L88").Value = 1.29
$ws.Range("M87").Value = '29/10/2023 14:26'
$ws.Range("M88").Value = '29/10/2023 14:23'
$ws.Range("N87").Value = 3.88
$ws.Range("N88").Value = 4.9
$ws.Range("O87").Value = '22/10/2023 14:42'
$ws.Range("O88").Value = '22/10/2023 12:43'
$ws.Range("P87").Value = 3.81
$ws.Range("P88").Value = 6.33
$ws.Range("Q87").Value = '29/10/2023 14:24'
$ws.Range("Q88").Value = '29/10/2023 14:29'
$ws.Range("R87").Value = 2.25
$ws.Range("R88").Value = 5.21
$ws.Range("S87").Value = '22/10/2023 14:42'
$ws.Range("S88").Value = '22/10/2023 12:43'
$ws.Range("T87").Value = 2.61
$ws.Range("T88").Value = 9.79
$ws.Range("U87").Value = '29/10/2023 14:21'
$ws.Range("U88").Value = '29/10/2023 14:29'
$ws.Range("V87").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/fc-volendam-excelsior/Onc5ROIp/'
$ws.Range("V88").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-ajax/I5KxW4AT/'

# swap row 91 and 92
$ws.Range("F91").Value = 'Waalwijk'
$ws.Range("F92").Value = 'Excelsior'
$ws.Range("G91").Value = 1
$ws.Range("G92").Value = 1
$ws.Range("H91").Value = 'Feyenoord'
$ws.Range("H92").Value = 'AZ Alkmaar'
$ws.Range("I91").Value = 2
$ws.Range("I92").Value = 1
$ws.Range("J91").Value = 8.58
$ws.Range("J92").Value = 5.11
$ws.Range("K91").Value = '29/10/2023 12:42'
$ws.Range("K92").Value = '29/10/2023 17:13'
$ws.Range("L91").Value = 14.74
$ws.Range("L92").Value = 9.11
$ws.Range("M91").Value = '04/11/2023 18:40'
$ws.Range("M92").Value = '04/11/2023 18:43'
$ws.Range("N91").Value = 5.63
$ws.Range("N92").Value = 4.63
$ws.Range("O91").Value = '29/10/2023 12:42'
$ws.Range("O92").Value = '29/10/2023 17:13'
$ws.Range("P91").Value = 8.98
$ws.Range("P92").Value = 5.56
$ws.Range("Q91").Value = '04/11/2023 18:40'
$ws.Range("Q92").Value = '04/11/2023 18:43'
$ws.Range("R91").Value = 1.33
$ws.Range("R92").Value = 1.6
$ws.Range("S91").Value = '29/10/2023 12:42'
$ws.Range("S92").Value = '29/10/2023 17:13'
$ws.Range("T91").Value = 1.17
$ws.Range("T92").Value = 1.35
$ws.Range("U91").Value = '04/11/2023 18:34'
$ws.Range("U92").Value = '04/11/2023 18:35'
$ws.Range("V91").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-feyenoord/C8w0m6u4/'
$ws.Range("V92").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-az-alkmaar/UqlelnQc/'

# swap row 97 and 98
$ws.Range("F97").Value = 'Ajax'
$ws.Range("F98").Value = 'Sparta Rotterdam'
$ws.Range("G97").Value = 4
$ws.Range("G98").Value = 1
$ws.Range("H97").Value = 'Heerenveen'
$ws.Range("H98").Value = 'Almere City'
$ws.Range("I97").Value = 1
$ws.Range("I98").Value = 2
$ws.Range("J97").Value = 1.62
$ws.Range("J98").Value = 2.01
$ws.Range("K97").Value = '02/11/2023 20:12'
$ws.Range("K98").Value = '28/10/2023 20:13'
$ws.Range("L97").Value = 1.41
$ws.Range("L98").Value = 1.75
$ws.Range("M97").Value = '05/11/2023 16:44'
$ws.Range("M98").Value = '05/11/2023 16:36'
$ws.Range("N97").Value = 4.85
$ws.Range("N98").Value = 3.74
$ws.Range("O97").Value = '02/11/2023 20:12'
$ws.Range("O98").Value = '28/10/2023 20:13'
$ws.Range("P97").Value = 5.6
$ws.Range("P98").Value = 3.98
$ws.Range("Q97").Value = '05/11/2023 16:44'
$ws.Range("Q98").Value = '05/11/2023 16:36'
$ws.Range("R97").Value = 4.55
$ws.Range("R98").Value = 3.73
$ws.Range("S97").Value = '02/11/2023 20:12'
$ws.Range("S98").Value = '28/10/2023 20:13'
$ws.Range("T97").Value = 7.02
$ws.Range("T98").Value = 4.79
$ws.Range("U97").Value = '05/11/2023 16:44'
$ws.Range("U98").Value = '05/11/2023 16:36'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/ajax-heerenveen/jgkikSAi/'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-almere-city/z7ruv22j/'
# --- Append 3 new rows (100,101,102) with formatting copied from row 99 ---
$ws.Range("A99:V99").Copy($ws.Range("A100:V100"))
$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 'netherlands'
$ws.Range("C100").Value = 'eredivisie'
$ws.Range("D100").Value = '2023-2024'
$ws.Range("E100").Value = 45241.78125
$ws.Range("F100").Value = 'Waalwijk'
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 'G.A. Eagles'
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 2.53
$ws.Range("K100").Value = '04/11/2023 21:12'
$ws.Range("L100").Value = 3.22
$ws.Range("M100").Value = '11/11/2023 18:37'
$ws.Range("N100").Value = 3.55
$ws.Range("O100").Value = '04/11/2023 21:12'
$ws.Range("P100").Value = 3.8
$ws.Range("Q100").Value = '11/11/2023 18:37'
$ws.Range("R100").Value = 2.76
$ws.Range("S100").Value = '04/11/2023 21:12'
$ws.Range("T100").Value = 2.2
$ws.Range("U100").Value = '11/11/2023 18:37'
$ws.Range("V100").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-g-a-eagles/08GXqhj9/'

$ws.Range("A99:V99").Copy($ws.Range("A101:V101"))
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 'netherlands'
$ws.Range("C101").Value = 'eredivisie'
$ws.Range("D101").Value = '2023-2024'
$ws.Range("E101").Value = 45241.83333333334
$ws.Range("F101").Value = 'Twente'
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 'Nijmegen'
$ws.Range("I101").Value = 3
$ws.Range("J101").Value = 1.35
$ws.Range("K101").Value = '05/11/2023 14:42'
$ws.Range("L101").Value = 1.35
$ws.Range("M101").Value = '11/11/2023 19:44'
$ws.Range("N101").Value = 5.56
$ws.Range("O101").Value = '05/11/2023 14:42'
$ws.Range("P101").Value = 5.45
$ws.Range("Q101").Value = '11/11/2023 19:59'
$ws.Range("R101").Value = 8.07
$ws.Range("S101").Value = '05/11/2023 14:42'
$ws.Range("T101").Value = 9.37
$ws.Range("U101").Value = '11/11/2023 19:59'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-nijmegen/Isn9Q3kl/'

$ws.Range("A99:V99").Copy($ws.Range("A102:V102"))
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 'netherlands'
$ws.Range("C102").Value = 'eredivisie'
$ws.Range("D102").Value = '2023-2024'
$ws.Range("E102").Value = 45241.875
$ws.Range("F102").Value = 'Vitesse'
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 'Heerenveen'
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 2.21
$ws.Range("K102").Value = '05/11/2023 17:12'
$ws.Range("L102").Value = 2.39
$ws.Range("M102").Value = '11/11/2023 20:52'
$ws.Range("N102").Value = 3.69
$ws.Range("O102").Value = '05/11/2023 17:12'
$ws.Range("P102").Value = 3.33
$ws.Range("Q102").Value = '11/11/2023 20:52'
$ws.Range("R102").Value = 3.16
$ws.Range("S102").Value = '05/11/2023 17:12'
$ws.Range("T102").Value = 3.23
$ws.Range("U102").Value = '11/11/2023 20:58'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/vitesse-heerenveen/hrDPoELd/'
